# 200m Bryst_statistics.xlsx
#
# The swimmer "Bjarne Forfot" (time 3.27,19 / 195 points / 04.12.2016 /
# Stjørdal / 25m pool) had been recorded with the wrong gender: his row
# lived on the Female_25m sheet (at row 128) instead of Male_25m.
#
# Fix: remove the misfiled row from Female_25m and insert the same data,
# with the correct "Male" gender, into Male_25m at its sorted position
# (row 112, by points descending).

$wb = $excel.ActiveWorkbook

# --- Male_25m: insert the corrected row for Bjarne Forfot -------------
$maleSheet = $wb.Worksheets.Item("Male_25m")
$maleSheet.Rows("112:112").Insert()

# Make sure the "Dato" column keeps storing the date as literal text
# (matching every other row) instead of Excel auto-coercing the
# dd.mm.yyyy-looking string into a real date serial number.
$maleSheet.Cells.Item(112, 4).NumberFormat = "@"

$maleSheet.Cells.Item(112, 1).Value = "Bjarne Forfot"
$maleSheet.Cells.Item(112, 2).Value = "3.27,19"
$maleSheet.Cells.Item(112, 3).Value = 195
$maleSheet.Cells.Item(112, 4).Value = "04.12.2016"
$maleSheet.Cells.Item(112, 5).Value = "Stjørdal"
$maleSheet.Cells.Item(112, 6).Value = "25m"
$maleSheet.Cells.Item(112, 7).Value = "Male"

# Drop the temporary text format again so the new row is unstyled, just
# like every other data row on this sheet.
$maleSheet.Range("A112:G112").ClearFormats()

# --- Female_25m: remove the misfiled "Bjarne Forfot" row --------------
$femaleSheet = $wb.Worksheets.Item("Female_25m")
$femaleSheet.Rows("128:128").Delete()
